$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.276.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "'1.783.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'340.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "'0.3957"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.53%  "

$ws.Range("D8").Value = "'0.3459"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.13%  "

$ws.Range("D9").Value = "'47.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("D10").Value = "'1.196"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.38%  "

$ws.Range("D11").Value = "'0.07464"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.60%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "'21.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.60%  "

$ws.Range("D14").Value = "'6.487"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.07%  "

$ws.Range("D15").Value = "'1.782.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.31%  "

$ws.Range("D16").Value = "'7.105"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("E17").Value = "  -2.94%  "

$ws.Range("D18").Value = "'0.06692"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").Value = "'84.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.09%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").Value = "'17.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").Value = "'6.511"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").Value = "'27.286.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").Value = "'12.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.77%  "

$ws.Range("D25").Value = "'2.381"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.13%  "

$ws.Range("D26").Value = "'1.470"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").Value = "'21.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.10%  "

$ws.Range("E28").Value = "  -6.92%  "

$ws.Range("D29").Value = "'157.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.01%  "

$ws.Range("D30").Value = "'1.985.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("D31").Value = "'136.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("D32").Value = "'4.029"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("D33").Value = "'5.964"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.91%  "

$ws.Range("D34").Value = "'0.08826"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("E35").Value = "  -6.69%  "

$ws.Range("D36").Value = "'0.02440"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.31%  "

$ws.Range("D37").Value = "'1.621"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.49%  "

$ws.Range("D38").Value = "'5.410"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.01%  "

$ws.Range("D39").Value = "'0.06460"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.93%  "

$ws.Range("D40").Value = "'0.6822"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.92%  "

$ws.Range("D41").Value = "'0.2207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "

$ws.Range("D42").Value = "'1.251"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.18%  "

$ws.Range("D43").Value = "'8.391"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.16%  "

$ws.Range("D44").Value = "'14.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("D46").Value = "'0.6390"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("D47").Value = "'3.880"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.60%  "

$ws.Range("D48").Value = "'2.136"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.57%  "

$ws.Range("D49").Value = "'132.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").Value = "'0.07139"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.34%  "

$ws.Range("D51").Value = "'79.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.44%  "
